# Auto-update draw results: append the 2025-11-16 Pick 3 draw as a new row
# (row 61) at the bottom of the "Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A61:E61")

# Temporarily force the new cells to Text format so that values which look
# like dates/numbers (the draw date "2025-11-16" and the "YYMMDD" phase code
# "251116") are stored as literal strings rather than being auto-coerced by
# Excel's smart input parsing into a real date serial / number - matching how
# every other row in this sheet is stored.
$rng.NumberFormat = "@"

$ws.Range("A61").Value = "2025-11-16"
$ws.Range("B61").Value = "Pick 3"
$ws.Range("C61").Value = "251116"
$ws.Range("D61").Value = "4-2-9"
$ws.Range("E61").Value = "2025-11-16T21:37:03.249+04:00"

# Remove the temporary explicit format again so the new row's cells end up
# with the same "no explicit style" footprint as all the existing rows.
$rng.ClearFormats()
